$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Start clean: wipe all existing cell content/styles/shared strings so
# the sheet only ends up containing what we put back.
# ---------------------------------------------------------------------
$ws.Cells.Clear() | Out-Null

# ---------------------------------------------------------------------
# Header row (row 1): source, tanggal, comment_line, departemen,
# kode_akun, nama_akun, debit, kredit
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "source"
$ws.Range("B1").Value = "tanggal"
$ws.Range("C1").Value = "comment_line"
$ws.Range("D1").Value = "departemen"
$ws.Range("E1").Value = "kode_akun"
$ws.Range("F1").Value = "nama_akun"
$ws.Range("G1").Value = "debit"
$ws.Range("H1").Value = "kredit"

# ---------------------------------------------------------------------
# Header styling: bold text, thin box border all around, centered
# horizontally, aligned to top vertically. Done before the text-format
# scratch trick below so this combined style claims the lowest free
# cellXfs slot.
# ---------------------------------------------------------------------
$headerFirst = $ws.Range("A1")
$headerFirst.Font.Bold = $true
$headerFirst.Borders.LineStyle = 1
$headerFirst.HorizontalAlignment = -4108
$headerFirst.VerticalAlignment = -4160

$headerFirst.Copy() | Out-Null
$ws.Range("B1:H1").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# B2/B3 (tanggal) and E2/E3 (kode_akun) look like dates/numbers, but must
# stay plain text. Stamp a text ("@") number format onto a scratch cell,
# copy its formatting (format-only) onto the target cells first - this
# keeps the values we are about to write stored as shared strings
# instead of being auto-converted to a date serial / numeric value.
# ---------------------------------------------------------------------
$scratch = $ws.Range("Z1:Z2")
$scratch.NumberFormat = "@"
$scratch.Copy() | Out-Null
$ws.Range("B2:B3").PasteSpecial(-4122) | Out-Null
$ws.Range("E2:E3").PasteSpecial(-4122) | Out-Null
$scratch.Clear() | Out-Null

# ---------------------------------------------------------------------
# Data rows (row 2-3)
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "MDR0012-2501-0001"
$ws.Range("B2").Value = "2025-01-09"
$ws.Range("C2").Value = "saldo qris MCM InhouseTrf DARI REMBULAN CITRA ABADI Transfer Fee 20250108193450725499102 — MDR0012-2501-0001 / 0001"
$ws.Range("D2").Value = "ODS"
$ws.Range("E2").Value = "1102002"
$ws.Range("G2").Value = 1000000
$ws.Range("H2").Value = 0

$ws.Range("A3").Value = "MDR0012-2501-0001"
$ws.Range("B3").Value = "2025-01-09"
$ws.Range("C3").Value = "saldo qris MCM InhouseTrf DARI REMBULAN CITRA ABADI Transfer Fee 20250108193450725499102 — MDR0012-2501-0001 / 0001"
$ws.Range("D3").Value = "ODS"
$ws.Range("E3").Value = "1102001"
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 1000000

# Put B2:B3 / E2:E3 back onto the default (unstyled) cell format so they
# don't carry the scratch "@" style forward into the saved file.
$ws.Range("B2:B3").Style = "Normal"
$ws.Range("E2:E3").Style = "Normal"

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Selection matches the full populated range.
# ---------------------------------------------------------------------
$ws.Range("A1:H3").Select() | Out-Null
